{"js": "// The supplied unified diff for this commit touches ONLY the package's\n// internal customXml parts:\n//   - customXml/item1.xml .. item4.xml            (content of the 4 parts\n//     got cyclically re-numbered: item1<->item4, item2<->item3)\n//   - customXml/itemProps1.xml .. itemProps4.xml   (re-numbered the same\n//     way, to stay paired with their item*.xml via the unchanged\n//     customXml/_rels/item*.xml.rels)\n//\n// Those four parts are SharePoint/content-type/bibliography plumbing\n// (FormTemplates, a ct:contentTypeSchema, a b:Sources bibliography stub,\n// and a p:properties/_activity stub) \u2014 none of it is visible document\n// content, and the four XML payloads themselves are byte-identical\n// before/after the commit; only their item-number slot changed. That is\n// an artifact of Word/SharePoint re-saving the file (the commit message\n// is about renaming an unrelated repository folder,\n// \"Process_Phase\" -> \"Process&Analyze_Phase\", that this file does not\n// even live in \u2014 it lives under \"Ask_Phase/\").\n//\n// The diff contains no hunks for word/document.xml, styles.xml,\n// numbering.xml, headers/footers, etc., so the document's visible\n// text/formatting/structure is unchanged. The Word JavaScript API only\n// exposes the document's content model (body, paragraphs, tables,\n// ranges, ...), not raw OPC package/custom-XML-part renumbering, so\n// there is nothing reachable through `context.document` that needs to\n// change to match this diff. We simply touch the body (load/sync) to\n// confirm the document is reachable and leave its content untouched.\n\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The supplied unified diff for this commit touches ONLY the package's\n# internal customXml parts:\n#   - customXml/item1.xml .. item4.xml            (content of the 4 parts\n#     got cyclically re-numbered: item1<->item4, item2<->item3)\n#   - customXml/itemProps1.xml .. itemProps4.xml   (re-numbered the same\n#     way, to stay paired with their item*.xml via the unchanged\n#     customXml/_rels/item*.xml.rels)\n#\n# Those four parts are SharePoint/content-type/bibliography plumbing\n# (FormTemplates, a ct:contentTypeSchema, a b:Sources bibliography stub,\n# and a p:properties/_activity stub) - none of it is visible document\n# content, and the four XML payloads themselves are byte-identical\n# before/after the commit; only their item-number slot changed. That is\n# an artifact of Word/SharePoint re-saving the file (the commit message\n# is about renaming an unrelated repository folder,\n# \"Process_Phase\" -> \"Process&Analyze_Phase\", that this file does not\n# even live in - it lives under \"Ask_Phase/\").\n#\n# The diff contains no hunks for word/document.xml, styles.xml,\n# numbering.xml, headers/footers, etc., so the document's visible\n# text/formatting/structure is unchanged. The Word COM object model's\n# CustomXMLParts collection is for custom XML data parts, not for\n# reordering/renumbering the underlying OPC package parts, and there is\n# no visible-content change to make here. We simply touch the document\n# (read its content) to confirm it is reachable and leave it untouched.\n\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
